$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2023" column (T) that mirrors the formatting of the existing
# 2022 column (S), then fill in the new figures for each indicator row.
$ws.Range("S4:S8").Copy()
$ws.Range("T4:T8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 43.1
$ws.Range("T6").Value = 19.7
$ws.Range("T7").Value = 7.8
$ws.Range("T8").Value = 15.6

# Row heights were tightened slightly when the new column was added.
$ws.Rows.Item(1).RowHeight = 57
$ws.Rows.Item(4).RowHeight = 16.5

# Row 9 has no data in the new column, but its row element still needs to
# report the widened column span; touching the row (without altering any
# cell) is enough to refresh it.
$ws.Rows.Item(9).EntireRow.Select()

# Leave the selection on a sensible, stable cell before saving.
$ws.Range("A1").Select()
